# report_info.xlsx edit script
# Commit message: "AG principal added to NPIV table"

$wb = $excel.ActiveWorkbook

# --- Sheet "report" (sheet1) ---
$ws1 = $wb.Worksheets.Item("report")

# Move the "value" header cell from D1 to C1 (keep style/value, clear old cell)
$ws1.Range("D1").Copy($ws1.Range("C1"))
$ws1.Range("D1").Clear()

# Update the saved selection / active cell for this sheet
$ws1.Range("C5").Select() | Out-Null

# --- Sheet "service_tables" (sheet2) ---
$ws2 = $wb.Worksheets.Item("service_tables")
$ws2.Activate() | Out-Null

# peerzone / peerzone_effective rows - disable export/force flags
$ws2.Range("G30").Value = 0
$ws2.Range("H30").Value = 0
$ws2.Range("G31").Value = 0

# NPIV row - enable export flag (AG principal added to NPIV table)
$ws2.Range("G62").Value = 1

# Zoning configuration data-analysis rows - enable export flags
$ws2.Range("H70").Value = 1
$ws2.Range("G71").Value = 1
$ws2.Range("G72").Value = 1
$ws2.Range("G73").Value = 1
$ws2.Range("G74").Value = 1
$ws2.Range("G75").Value = 1
$ws2.Range("G76").Value = 1
$ws2.Range("G77").Value = 1
$ws2.Range("G78").Value = 1
$ws2.Range("G81").Value = 1

# Add an AutoFilter over the full table range
$ws2.Range("A1:J81").AutoFilter() | Out-Null

# Update the view (scroll position + selection) for this sheet
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 55
$win2.ScrollColumn = 3
$ws2.Range("E75").Select() | Out-Null
